# Apply the edits described by the diff.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date and FHIR Version values ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$wsMeta.Range("B15").Value = "4.0.1"

# --- Elements sheet: update constraint / type / reference text ---
$wsElem = $wb.Worksheets.Item("Elements")

$eleShort = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension row (row 2) constraint text loses the "unless an empty Parameters resource ... " clause
$wsElem.Range("AJ2").Value = $eleShort

# Extension.id row (row 3): Type changes from "id" to "string"
$wsElem.Range("K3").Value = "string" + [char]10

# Extension.extension row (row 4) constraint text matches the shortened text too
$wsElem.Range("AJ4").Value = $eleShort

# Extension.value[x] row (row 6): Definition text now references R4 instead of R4B
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
